$d = $word.ActiveDocument

# The pkg:package wrapper namespace used by Range.InsertXML to splice raw
# WordprocessingML back into word/document.xml.
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# Edit 1: address block paragraph "01010, м. Київ, вул. Московська, 8"
#   - drop the stray <w:lang w:val="ru-RU"/> on the ", " run
#   - split "м. Київ, вул. Московська, 8" into its own two paragraphs
#     ("м. Київ, " stays put, "вул. Московська, 8" becomes a new line)
# ---------------------------------------------------------------------
$addr = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Московська*") {
        $addr = $p.Range
    }
}
if ($addr -ne $null) {
    $addrXml = $pkgOpen + '<w:body>' +
        '<w:p><w:pPr><w:ind w:left="5760"/><w:rPr><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>01010</w:t></w:r>' +
        '<w:r><w:rPr><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' +
        '<w:r><w:rPr><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">м. Київ, </w:t></w:r>' +
        '</w:p>' +
        '<w:p><w:pPr><w:ind w:left="5760"/><w:rPr><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>вул. Московська, 8</w:t></w:r>' +
        '</w:p>' +
        '</w:body>' + $pkgClose
    $addr.InsertXML($addrXml)
}

# ---------------------------------------------------------------------
# Edit 2: "Вказане рішення було подано ... 16 серпня 2018 року ..."
#   - replace the literal date "16 серпня 2018" with a "{5}" placeholder
#     run (kept separate from the surrounding text runs)
# ---------------------------------------------------------------------
$decision = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Вказане рішення*") {
        $decision = $p.Range
    }
}
if ($decision -ne $null) {
    $decisionXml = $pkgOpen + '<w:body>' +
        '<w:p><w:pPr><w:pStyle w:val="ab"/><w:ind w:firstLine="567"/><w:jc w:val="both"/>' +
        '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:bCs/><w:lang w:val="uk-UA"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:bCs/><w:lang w:val="uk-UA"/></w:rPr>' +
        '<w:t xml:space="preserve">Вказане рішення було подано державному реєстратору </w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:bCs/><w:lang w:val="uk-UA"/></w:rPr>' +
        '<w:t>{5}</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:bCs/><w:lang w:val="uk-UA"/></w:rPr>' +
        '<w:t xml:space="preserve"> року для проведення реєстраційної дії «Внесення рішення засновників (учасників) юридичної особи або уповноваженого ними органу щодо припинення юридичної особи». </w:t></w:r>' +
        '</w:p>' +
        '</w:body>' + $pkgClose
    $decision.InsertXML($decisionXml)
}
